# "Colocando header nos gráficos"
# Adds header labels in column A (row 1) for the scenario sheets, normalizes
# accented Portuguese text that was previously stored without diacritics,
# removes the bold/bordered header style from the data-label cells in column A
# (rows 2-12), drops the now-unused "Teto" row from the emissions sheet, and
# refreshes the cost sheet header/values.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheets 1-4 share the same layout:
#   A1 (new header) | B1..E1 = years
#   A2..A12 = source/technology labels, B..E = values
# ---------------------------------------------------------------------------
$sheetIndexes = @(1, 2, 3, 4)

$labels = @{
    2  = "Hidro"
    3  = "Gás Natural"
    4  = "Carvão"
    5  = "Nuclear"
    6  = "Óleos Comb"
    7  = "Biomassa"
    8  = "Eólica"
    9  = "Solar"
    10 = "Outros"
    11 = "Pot. Compl."
    12 = "GD"
}

foreach ($idx in $sheetIndexes) {
    $ws = $wb.Worksheets.Item($idx)

    # Add the new header cell in A1, copying the existing bold/border style
    # used by the rest of row 1 (B1:E1) so it matches visually.
    $ws.Range("B1").Copy($ws.Range("A1"))
    $ws.Range("A1").Value = "Fonte/Tecnologia"

    # Update the row labels (fixing missing accents) and strip the header
    # style that had been applied to them (they are plain data labels now).
    foreach ($row in $labels.Keys) {
        $cell = $ws.Cells.Item($row, 1)
        $cell.Value = $labels[$row]
        $cell.ClearFormats()
    }
}

# ---------------------------------------------------------------------------
# Sheet 5 - "Emissoes Totais (MtCO2eq)"
#   Add header "Período" to A1, fix accents on A2/A3 labels, remove the
#   unused "Teto" row (row 4).
# ---------------------------------------------------------------------------
$ws5 = $wb.Worksheets.Item(5)

$ws5.Range("B1").Copy($ws5.Range("A1"))
$ws5.Range("A1").Value = "Período"

$ws5.Cells.Item(2, 1).Value = "P.Médio"
$ws5.Cells.Item(2, 1).ClearFormats()

$ws5.Cells.Item(3, 1).Value = "P.Crítico"
$ws5.Cells.Item(3, 1).ClearFormats()

$ws5.Rows.Item(4).Delete()

# ---------------------------------------------------------------------------
# Sheet 6 - "Custo Total (bilhões de R$)"
#   Add header "Tipo Expansão" to A1, change B1 from "Custo" to the (text)
#   label "2015", fix accents on A2/A3 labels, update B2/B3 values.
# ---------------------------------------------------------------------------
$ws6 = $wb.Worksheets.Item(6)

$ws6.Range("B1").Copy($ws6.Range("A1"))
$ws6.Range("A1").Value = "Tipo Expansão"

# Force B1 to hold the text "2015" (not the number 2015), matching how the
# year headers are stored as text on the other sheets. A plain
# Range.Value = "2015" would be auto-converted to a number by Excel, so the
# text value is first produced in a scratch cell (formatted as Text) and
# only its value is pasted into B1, leaving B1's original bold/border/center
# header style (copied from A1's old style, still s="1") untouched.
$scratch = $ws6.Cells.Item(100, 100)
$scratch.NumberFormat = "@"
$scratch.Value = "2015"
$scratch.Copy()
$ws6.Range("B1").PasteSpecial(-4163)
$scratch.Clear()

$ws6.Cells.Item(2, 1).Value = "Expansão Centralizada"
$ws6.Cells.Item(2, 1).ClearFormats()
$ws6.Cells.Item(2, 2).Value = 588

$ws6.Cells.Item(3, 1).Value = "Expansão por GD"
$ws6.Cells.Item(3, 1).ClearFormats()
$ws6.Cells.Item(3, 2).Value = 99
